$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Globo"
$ws.Range("B7").Value = "RJ TV 1"
$ws.Range("C7").Value = "Social"
$ws.Range("D7").Value = "2025-04-01T13:31"
$ws.Range("E7").Value = "Positivo"
$ws.Range("F7").Value = "Curso Ferroport + Firjan - Campos e SJB - Técnicas em manutenção e eletromecânica. Preferencialmente, para mulheres e pessoas negras"
